# Applies the "📊 Excel mis à jour automatiquement" update:
#   - "Prix Spot" sheet: append a new day column (CM) "12-sep" with its 24 hourly values
#   - "Gaz" sheet: append a new row (88) for 2025-09-10
#   - "CO2" sheet: append a new row (88) for 2025-09-10

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Prix Spot": new column CM (day "12-sep"), rows 1 (header) .. 25
# ---------------------------------------------------------------------------
$wsPrix = $wb.Worksheets.Item("Prix Spot")

# Header cell, formatted like the preceding day header (CL1: bold/centered/bordered)
$wsPrix.Cells.Item(1, 91).Value = "12-sep"
$wsPrix.Range("CL1").Copy()
$wsPrix.Range("CM1").PasteSpecial(-4122)

$prixValues = @{
    2  = 0
    3  = 0.01
    4  = 0.44
    5  = 0.01
    6  = 0.01
    7  = 3
    8  = 11.58
    9  = 21.04
    10 = 41.82
    11 = 26.68
    12 = 1.85
    13 = 0
    14 = 1.77
    15 = 0
    16 = 0
    17 = 0
    18 = 3.52
    19 = 19.68
    20 = 39.92
    21 = 61.68
    22 = 66.34
    23 = 46
    24 = 49.97
    25 = 41.2
}

foreach ($row in $prixValues.Keys) {
    $wsPrix.Cells.Item($row, 91).Value = $prixValues[$row]
}

# ---------------------------------------------------------------------------
# Sheet "Gaz": new row 88 -> 2025-09-10 / 32.6
# ---------------------------------------------------------------------------
$wsGaz = $wb.Worksheets.Item("Gaz")

# Format the date cell as text first so the "yyyy-mm-dd" string isn't
# auto-converted into a date serial number, then drop the number format
# again so the cell keeps the default (unstyled) look of its neighbours.
$wsGaz.Cells.Item(88, 1).NumberFormat = "@"
$wsGaz.Cells.Item(88, 1).Value = "2025-09-10"
$wsGaz.Cells.Item(88, 1).ClearFormats()
$wsGaz.Cells.Item(88, 2).Value = 32.6

# ---------------------------------------------------------------------------
# Sheet "CO2": new row 88 -> 2025-09-10 / 76.2
# ---------------------------------------------------------------------------
$wsCo2 = $wb.Worksheets.Item("CO2")

$wsCo2.Cells.Item(88, 1).NumberFormat = "@"
$wsCo2.Cells.Item(88, 1).Value = "2025-09-10"
$wsCo2.Cells.Item(88, 1).ClearFormats()
$wsCo2.Cells.Item(88, 2).Value = 76.2
